# Update the "想去人数" (want-to-go count) figures in F column across the
# three affected sheets, matching a fresh scrape of the source data.
# Sheet order (per workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value = 2800
$wsExpo.Range("F4").Value = 20450
$wsExpo.Range("F6").Value = 2556
$wsExpo.Range("F7").Value = 782
$wsExpo.Range("F10").Value = 731
$wsExpo.Range("F16").Value = 494

# --- Sheet 3: 本地生活 (Local life) ---
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Range("F2").Value = 6076
$wsLocal.Range("F4").Value = 648
$wsLocal.Range("F5").Value = 1378
$wsLocal.Range("F6").Value = 38

# --- Sheet 4: 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 6076
$wsAll.Range("F4").Value = 648
$wsAll.Range("F5").Value = 1378
$wsAll.Range("F6").Value = 2800
$wsAll.Range("F8").Value = 20450
$wsAll.Range("F14").Value = 2556
$wsAll.Range("F15").Value = 782
$wsAll.Range("F17").Value = 38
$wsAll.Range("F20").Value = 731
$wsAll.Range("F31").Value = 494
